$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Skills" column before column D (Location), shifting Location/experience/salary/Description right
$ws.Columns.Item(4).Insert()
$ws.Range("D1").Value = "Skills"

# Update data rows 2-19 with new values (Title/B unchanged column letter, Companies/C unchanged, new Skills/D, Location/E, experience/F, salary/G, Description/H)
$ws.Cells.Item(2, 2).Value = "Computer Science Engineering- Professor"
$ws.Cells.Item(2, 3).Value = "Sharda University3.6248  Reviews"
$ws.Cells.Item(2, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, advisory, computer"
$ws.Cells.Item(2, 5).Value = "Greater Noida"
$ws.Cells.Item(2, 6).Value = "10-12 Yrs"
$ws.Cells.Item(2, 7).Value = "Not disclosed"
$ws.Cells.Item(2, 8).Value = "experience as a Professor including 5 years of administrative experience in a reputed U..."

$ws.Cells.Item(3, 2).Value = "Computer Science & Informatics- Professor"
$ws.Cells.Item(3, 3).Value = "Maharshi Dayanand University (MDU)4.819  Reviews"
$ws.Cells.Item(3, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, Education, Computer"
$ws.Cells.Item(3, 5).Value = "Jaipur"
$ws.Cells.Item(3, 6).Value = "7-10 Yrs"
$ws.Cells.Item(3, 7).Value = "Not disclosed"
$ws.Cells.Item(3, 8).Value = "Work with students who are studying for a degree or a certificate or certification or a..."

$ws.Cells.Item(4, 2).Value = "Computer Science & Informatics- Professor"
$ws.Cells.Item(4, 3).Value = "Maharshi Dayanand University (MDU)4.819  Reviews"
$ws.Cells.Item(4, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, Education, Computer"
$ws.Cells.Item(4, 5).Value = "Jaipur"
$ws.Cells.Item(4, 6).Value = "8-11 Yrs"
$ws.Cells.Item(4, 7).Value = "Not disclosed"
$ws.Cells.Item(4, 8).Value = "Work with students who are studying for a degree or a certificate or certification or a..."

$ws.Cells.Item(5, 2).Value = "Sr . Software Developer"
$ws.Cells.Item(5, 3).Value = "Newstar Infotech"
$ws.Cells.Item(5, 4).Value = "Computer science, Web technologies, SQL database, Programming, Database, Application software, IOS, DBMS"
$ws.Cells.Item(5, 5).Value = "Ahmedabad"
$ws.Cells.Item(5, 6).Value = "0-1 Yrs"
$ws.Cells.Item(5, 7).Value = "Not disclosed"
$ws.Cells.Item(5, 8).Value = "Bachelor s degree or equivalent experience in Computer Science or related field Develop..."

$ws.Cells.Item(6, 2).Value = "Software Developer"
$ws.Cells.Item(6, 3).Value = "Hi Tech3.8135  Reviews"
$ws.Cells.Item(6, 4).Value = "C++, development, software, Development Manager, AutoCAD, VC++, CAD, Program Executive"
$ws.Cells.Item(6, 5).Value = "Chennai"
$ws.Cells.Item(6, 6).Value = "0-1 Yrs"
$ws.Cells.Item(6, 7).Value = "Not disclosed"
$ws.Cells.Item(6, 8).Value = "- Must be ambitious and have a desire to learn new skills"

$ws.Cells.Item(7, 2).Value = "Computer Science Faculty"
$ws.Cells.Item(7, 3).Value = "Amity University3.6966  Reviews"
$ws.Cells.Item(7, 4).Value = "Computer Science, Information Technology, Science, Technology, Computer"
$ws.Cells.Item(7, 5).Value = "Mohali, Punjab"
$ws.Cells.Item(7, 6).Value = "3-8 Yrs"
$ws.Cells.Item(7, 7).Value = "1.5-6 Lacs PA"
$ws.Cells.Item(7, 8).Value = "Teaching: Provide high-quality instruction at the undergraduate level, ensuring course ..."

$ws.Cells.Item(8, 2).Value = "Computer Science- TGT"
$ws.Cells.Item(8, 3).Value = "DPS Gandhinagar3.91694  Reviews"
$ws.Cells.Item(8, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, advisory, tgt"
$ws.Cells.Item(8, 5).Value = "Gandhinagar"
$ws.Cells.Item(8, 6).Value = "3-6 Yrs"
$ws.Cells.Item(8, 7).Value = "Not disclosed"
$ws.Cells.Item(8, 8).Value = "Job Description: ? Must possess a B.Ed. degree from recognized institution ? Schooling ..."

$ws.Cells.Item(9, 2).Value = "Computer Science - PGT"
$ws.Cells.Item(9, 3).Value = "DPS Patna3.91694  Reviews"
$ws.Cells.Item(9, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, advisory, computer"
$ws.Cells.Item(9, 5).Value = "Patna"
$ws.Cells.Item(9, 6).Value = "4-6 Yrs"
$ws.Cells.Item(9, 7).Value = "Not disclosed"
$ws.Cells.Item(9, 8).Value = "Job Description: ? Must possess a B.Ed. degree from recognized institution ? Schooling ..."

$ws.Cells.Item(10, 2).Value = "Computer Science Faculty"
$ws.Cells.Item(10, 3).Value = "Aditya P.U. College"
$ws.Cells.Item(10, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, Training, Science"
$ws.Cells.Item(10, 5).Value = "Bangalore/Bengaluru"
$ws.Cells.Item(10, 6).Value = "1-3 Yrs"
$ws.Cells.Item(10, 7).Value = "Not disclosed"
$ws.Cells.Item(10, 8).Value = "Key responsibilities of the job include: ? providing support to children in reading and..."

$ws.Cells.Item(11, 2).Value = "TGT Computer Science"
$ws.Cells.Item(11, 3).Value = "Zydus School For Excellence4.112  Reviews"
$ws.Cells.Item(11, 4).Value = "Professor, Tutor, English, Teachers, Education, Trainer, Lecturer, Teaching"
$ws.Cells.Item(11, 5).Value = "Godhavi"
$ws.Cells.Item(11, 6).Value = "3-5 Yrs"
$ws.Cells.Item(11, 7).Value = "Not disclosed"
$ws.Cells.Item(11, 8).Value = "Develop and maintain good Working habits and discipline in classroomsRequired TGT Compu..."

$ws.Cells.Item(12, 2).Value = "Professor - Computer Science"
$ws.Cells.Item(12, 3).Value = "Flame University4.621  Reviews"
$ws.Cells.Item(12, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, computer, computer science"
$ws.Cells.Item(12, 5).Value = "Bengaluru"
$ws.Cells.Item(12, 6).Value = "10-13 Yrs"
$ws.Cells.Item(12, 7).Value = "Not disclosed"
$ws.Cells.Item(12, 8).Value = "Work with students who are studying for a degree or a certificate or certification or a..."

$ws.Cells.Item(13, 2).Value = "Computer Science - Professor"
$ws.Cells.Item(13, 3).Value = "R. K. C. S. Educational Society"
$ws.Cells.Item(13, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, advisory, computer"
$ws.Cells.Item(13, 5).Value = "Firozabad"
$ws.Cells.Item(13, 6).Value = "10-12 Yrs"
$ws.Cells.Item(13, 7).Value = "Not disclosed"
$ws.Cells.Item(13, 8).Value = "A minimum of ten years of teaching experience in university / college, and / or experie..."

$ws.Cells.Item(14, 2).Value = "Computer Science - Professor"
$ws.Cells.Item(14, 3).Value = "R. K. C. S. Educational Society"
$ws.Cells.Item(14, 4).Value = "Counselor, Mentor, Trainer, Advisor, Educator, Teaching, advisory, computer"
$ws.Cells.Item(14, 5).Value = "New Delhi"
$ws.Cells.Item(14, 6).Value = "10-15 Yrs"
$ws.Cells.Item(14, 7).Value = "Not disclosed"
$ws.Cells.Item(14, 8).Value = "A minimum of ten years of teaching experience in university / college, and / or experie..."

$ws.Cells.Item(15, 2).Value = "Software Application Developer (OST)"
$ws.Cells.Item(15, 3).Value = "Logonb2b"
$ws.Cells.Item(15, 4).Value = "software, E-commerce, HTML, Research, Information technology, Joomla, Computer science, application"
$ws.Cells.Item(15, 5).Value = "Hyderabad"
$ws.Cells.Item(15, 6).Value = "0-2 Yrs"
$ws.Cells.Item(15, 7).Value = "Not disclosed"
$ws.Cells.Item(15, 8).Value = "Software Application Developer (OST) Skill: Good research, analytical, and Communicatio..."

$ws.Cells.Item(16, 2).Value = "Windows Application Software Developer"
$ws.Cells.Item(16, 3).Value = "Bluesurf Engineering Solutions"
$ws.Cells.Item(16, 4).Value = "Computer science, Agile scrum, Windows application, devops, Electronics, Telecommunication, Instrumentation, Application development"
$ws.Cells.Item(16, 5).Value = "Pune"
$ws.Cells.Item(16, 6).Value = "0-1 Yrs"
$ws.Cells.Item(16, 7).Value = "Not disclosed"
$ws.Cells.Item(16, 8).Value = "Bachelor s or master s degree in computer science / IT / Electronics / Electronics Tele..."

$ws.Cells.Item(17, 2).Value = "computer science Teacher AI, Machine Learning"
$ws.Cells.Item(17, 3).Value = "BDS Consultancy"
$ws.Cells.Item(17, 4).Value = "Computer Science, machine learning, ai, web development, cloud computing, Teaching, Computer, Machine"
$ws.Cells.Item(17, 5).Value = "Kolkata, West Bengal"
$ws.Cells.Item(17, 6).Value = "3-8 Yrs"
$ws.Cells.Item(17, 7).Value = "5-7.5 Lacs PA"
$ws.Cells.Item(17, 8).Value = "The candidate should have ME / M.Tech, Ph.D / Doctorate in Computer Science (Image Proc..."

$ws.Cells.Item(18, 2).Value = "Computer Science Faculty"
$ws.Cells.Item(18, 3).Value = "Vision Group of Colleges"
$ws.Cells.Item(18, 4).Value = "Computer Teaching, Lecturer Activities, Professor Activities, Education, Computer Science, Teaching, Information Technology, Computer"
$ws.Cells.Item(18, 5).Value = "Chittaurgarh, Rajasthan"
$ws.Cells.Item(18, 6).Value = "1-2 Yrs"
$ws.Cells.Item(18, 7).Value = "Not disclosed"
$ws.Cells.Item(18, 8).Value = "Deliver lectures, conduct practical sessions, and facilitate discussions for BCA (Bache..."

$ws.Cells.Item(19, 2).Value = "Computer Science Faculty"
$ws.Cells.Item(19, 3).Value = "Gnana Jyothi School3.03  Reviews"
$ws.Cells.Item(19, 4).Value = "Spoken English, Doctoral degree in computer science or a related field Minimum of five years of experience working in software development, Computer, English, Development, Computer science, Degree In Computer Science, Software"
$ws.Cells.Item(19, 5).Value = "Tumkur, Mysore/Mysuru, Bangalore/Bengaluru(Bagalur +2)"
$ws.Cells.Item(19, 6).Value = "3-5 Yrs"
$ws.Cells.Item(19, 7).Value = "3-4 Lacs PA"
$ws.Cells.Item(19, 8).Value = "Instruct students on how computers work, including the basic science and mathematics be..."

# Add two new rows (20, 21) for additional job listings, copying row 19 formatting for column A
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A19").Copy($ws.Range("A21"))

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Professor -  Computer Science & Engineering"
$ws.Cells.Item(20, 3).Value = "Sambhram Institute Of Technology2.916  Reviews"
$ws.Cells.Item(20, 4).Value = "Professor, Tutor, Teachers, Education, Trainer, Lecturer, Computer Science & Engineering, Teaching"
$ws.Cells.Item(20, 5).Value = "Bengaluru, Bangalore"
$ws.Cells.Item(20, 6).Value = "8-10 Yrs"
$ws.Cells.Item(20, 7).Value = "Not disclosed"
$ws.Cells.Item(20, 8).Value = "Required Professor for Computer Science & Engineering"

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Computer science Engineer - BHEL Jobs"
$ws.Cells.Item(21, 3).Value = "Right Step Consulting"
$ws.Cells.Item(21, 4).Value = "Software design, Programming, Software, TelecomDebugging, Application software, Linux kernelComp, Engineering, Computer science"
$ws.Cells.Item(21, 5).Value = "Noida"
$ws.Cells.Item(21, 6).Value = "1-5 Yrs"
$ws.Cells.Item(21, 7).Value = "Not disclosed"
$ws.Cells.Item(21, 8).Value = "Industry: Semiconductors / Electronics . Functional Area: IT Software - Telecom Softwar..."

